# Updated cryptos list with GitHub Actions
# Refreshes the latest price (column D) and 1h volume change (column E)
# snapshot values in the "cryptos" worksheet.
#
# Note: several Price values are plain decimal numbers (e.g. "4.46",
# "0.258") that Excel would otherwise auto-convert to numeric cells.
# The source data keeps these as literal text (matching the other
# dotted-thousands prices like "25.806.89"), so a leading apostrophe is
# used to force a text entry for those particular cells, exactly as a
# user typing in the grid would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.806.89"
$ws.Range("D3").Value = "1.636.75"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.258"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "1.862.26"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "1.636.83"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "25.826.06"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'4.46"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'192.67"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "'9.98"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "'142.41"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'0.905"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "1.131.39"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").Value = "'0.544"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "'100.55"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "1.771.68"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'55.33"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").Value = "'2.31"
$ws.Range("E51").Value = "  +2.58%  "
